$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEndTime = 45993.320254629631

# Rows 2-16: fixed "currently not charging" pinned terminals -- only refresh column D (latest charge-end-time)
$ws.Range("D2").Value = $newEndTime
$ws.Range("D3").Value = $newEndTime
$ws.Range("D4").Value = $newEndTime
$ws.Range("D5").Value = $newEndTime
$ws.Range("D6").Value = $newEndTime
$ws.Range("D7").Value = $newEndTime
$ws.Range("D8").Value = $newEndTime
$ws.Range("D9").Value = $newEndTime
$ws.Range("D10").Value = $newEndTime
$ws.Range("D11").Value = $newEndTime
$ws.Range("D12").Value = $newEndTime
$ws.Range("D13").Value = $newEndTime
$ws.Range("D14").Value = $newEndTime
$ws.Range("D15").Value = $newEndTime
$ws.Range("D16").Value = $newEndTime

# Row 17: pinned row with numeric C value, unchanged except D
$ws.Range("D17").Value = $newEndTime

# Rows 18-59: reshuffled set, sorted by column C ascending (re-fetched snapshot) -- rewrite A/B/C/D
$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "101号直流"
$ws.Range("C18").Value = 45987.552604166667
$ws.Range("D18").Value = $newEndTime

$ws.Range("A19").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B19").Value = "502号直流"
$ws.Range("C19").Value = 45989.209733796299
$ws.Range("D19").Value = $newEndTime

$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "603号直流"
$ws.Range("C20").Value = 45989.545405092591
$ws.Range("D20").Value = $newEndTime

$ws.Range("A21").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B21").Value = "406号直流"
$ws.Range("C21").Value = 45990.197141203702
$ws.Range("D21").Value = $newEndTime

$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "505号直流"
$ws.Range("C22").Value = 45991.02715277778
$ws.Range("D22").Value = $newEndTime

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "702号直流"
$ws.Range("C23").Value = 45991.033333333333
$ws.Range("D23").Value = $newEndTime

$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "703号直流"
$ws.Range("C24").Value = 45991.200173611112
$ws.Range("D24").Value = $newEndTime

$ws.Range("A25").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B25").Value = "301号直流"
$ws.Range("C25").Value = 45991.493321759262
$ws.Range("D25").Value = $newEndTime

$ws.Range("A26").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B26").Value = "102号直流"
$ws.Range("C26").Value = 45991.600381944445
$ws.Range("D26").Value = $newEndTime

$ws.Range("A27").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B27").Value = "903号直流"
$ws.Range("C27").Value = 45991.673692129632
$ws.Range("D27").Value = $newEndTime

$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "504号直流"
$ws.Range("C28").Value = 45992.035474537035
$ws.Range("D28").Value = $newEndTime

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "605号直流"
$ws.Range("C29").Value = 45992.047766203701
$ws.Range("D29").Value = $newEndTime

$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "A03号直流"
$ws.Range("C30").Value = 45992.202627314815
$ws.Range("D30").Value = $newEndTime

$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "602号直流"
$ws.Range("C31").Value = 45992.213622685187
$ws.Range("D31").Value = $newEndTime

$ws.Range("A32").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B32").Value = "503号直流"
$ws.Range("C32").Value = 45992.240590277775
$ws.Range("D32").Value = $newEndTime

$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "604号直流"
$ws.Range("C33").Value = 45992.250636574077
$ws.Range("D33").Value = $newEndTime

$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "402号直流"
$ws.Range("C34").Value = 45992.2731712963
$ws.Range("D34").Value = $newEndTime

$ws.Range("A35").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value = "006B号直流"
$ws.Range("C35").Value = 45992.277025462965
$ws.Range("D35").Value = $newEndTime

$ws.Range("A36").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B36").Value = "306号直流"
$ws.Range("C36").Value = 45992.404490740744
$ws.Range("D36").Value = $newEndTime

$ws.Range("A37").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B37").Value = "104号直流"
$ws.Range("C37").Value = 45992.429050925923
$ws.Range("D37").Value = $newEndTime

$ws.Range("A38").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B38").Value = "107号直流"
$ws.Range("C38").Value = 45992.519016203703
$ws.Range("D38").Value = $newEndTime

$ws.Range("A39").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B39").Value = "201号直流"
$ws.Range("C39").Value = 45992.520567129628
$ws.Range("D39").Value = $newEndTime

$ws.Range("A40").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B40").Value = "102号直流"
$ws.Range("C40").Value = 45992.523495370369
$ws.Range("D40").Value = $newEndTime

$ws.Range("A41").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B41").Value = "108号直流"
$ws.Range("C41").Value = 45992.526990740742
$ws.Range("D41").Value = $newEndTime

$ws.Range("A42").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value = "802号直流"
$ws.Range("C42").Value = 45992.537152777775
$ws.Range("D42").Value = $newEndTime

$ws.Range("A43").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B43").Value = "204号直流"
$ws.Range("C43").Value = 45992.54109953704
$ws.Range("D43").Value = $newEndTime

$ws.Range("A44").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B44").Value = "904号直流"
$ws.Range("C44").Value = 45992.545127314814
$ws.Range("D44").Value = $newEndTime

$ws.Range("A45").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B45").Value = "002A号直流"
$ws.Range("C45").Value = 45992.545925925922
$ws.Range("D45").Value = $newEndTime

$ws.Range("A46").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B46").Value = "805号直流"
$ws.Range("C46").Value = 45992.546932870369
$ws.Range("D46").Value = $newEndTime

$ws.Range("A47").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B47").Value = "705号直流"
$ws.Range("C47").Value = 45992.56077546296
$ws.Range("D47").Value = $newEndTime

$ws.Range("A48").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B48").Value = "106号直流"
$ws.Range("C48").Value = 45992.572858796295
$ws.Range("D48").Value = $newEndTime

$ws.Range("A49").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B49").Value = "205号直流"
$ws.Range("C49").Value = 45992.585555555554
$ws.Range("D49").Value = $newEndTime

$ws.Range("A50").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B50").Value = "210号直流"
$ws.Range("C50").Value = 45992.604259259257
$ws.Range("D50").Value = $newEndTime

$ws.Range("A51").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B51").Value = "A02号直流"
$ws.Range("C51").Value = 45992.63554398148
$ws.Range("D51").Value = $newEndTime

$ws.Range("A52").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B52").Value = "401号直流"
$ws.Range("C52").Value = 45992.650023148148
$ws.Range("D52").Value = $newEndTime

$ws.Range("A53").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B53").Value = "103号直流"
$ws.Range("C53").Value = 45992.655312499999
$ws.Range("D53").Value = $newEndTime

$ws.Range("A54").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B54").Value = "310号直流"
$ws.Range("C54").Value = 45992.656724537039
$ws.Range("D54").Value = $newEndTime

$ws.Range("A55").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B55").Value = "206号直流"
$ws.Range("C55").Value = 45992.663124999999
$ws.Range("D55").Value = $newEndTime

$ws.Range("A56").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B56").Value = "203号直流"
$ws.Range("C56").Value = 45992.676851851851
$ws.Range("D56").Value = $newEndTime

$ws.Range("A57").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B57").Value = "901号直流"
$ws.Range("C57").Value = 45992.678425925929
$ws.Range("D57").Value = $newEndTime

$ws.Range("A58").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B58").Value = "009A号直流"
$ws.Range("C58").Value = 45992.697662037041
$ws.Range("D58").Value = $newEndTime

$ws.Range("A59").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B59").Value = "105号直流"
$ws.Range("C59").Value = 45992.700428240743
$ws.Range("D59").Value = $newEndTime

# Update the active cell selection to match the author's final cursor position
$ws.Range("E13").Select()